$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "Dakota Myers" with "Franz Ferdinand" in cell B3
$ws.Range("B3").Value = "Franz Ferdinand"

# Update the active selection to B3 (was B17)
$ws.Range("B3").Select()
